# chore: update Sheets via scheduled runner
# Refreshes market-price-derived profit figures (columns H:N) on a handful
# of leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 166830.33
$ws.Range("I33").Value = 250160
$ws.Range("J33").Value = 171
$ws.Range("K33").Value = 250160
$ws.Range("L33").Value = 171
$ws.Range("M33").Value = -249931
$ws.Range("N33").Value = -629

$ws.Range("H38").Value = 1085
$ws.Range("I38").Value = 174
$ws.Range("J38").Value = 1996
$ws.Range("K38").Value = 522
$ws.Range("L38").Value = 5988
$ws.Range("M38").Value = -150
$ws.Range("N38").Value = -6732

$ws.Range("H40").Value = 35716610
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 38463890
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 38463890
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -38464240

$ws.Range("H55").Value = 305.23077
$ws.Range("I55").Value = 181.33333
$ws.Range("J55").Value = 411.42856
$ws.Range("K55").Value = 181.33333
$ws.Range("L55").Value = 411.42856
$ws.Range("M55").Value = 32.66667000000001
$ws.Range("N55").Value = -839.4285600000001

$ws.Range("H74").Value = 2040672.4
$ws.Range("I74").Value = 2317814
$ws.Range("J74").Value = 8300
$ws.Range("K74").Value = 2317814
$ws.Range("L74").Value = 8300
$ws.Range("M74").Value = -2316878
$ws.Range("N74").Value = -10172

$ws.Range("H77").Value = 2040672.4
$ws.Range("I77").Value = 2317814
$ws.Range("J77").Value = 8300
$ws.Range("K77").Value = 11589070
$ws.Range("L77").Value = 41500
$ws.Range("M77").Value = -11584390
$ws.Range("N77").Value = -50860

$ws.Range("H121").Value = 1599.8
$ws.Range("I121").Value = 800
$ws.Range("J121").Value = 2133
$ws.Range("K121").Value = 2400
$ws.Range("L121").Value = 6399
$ws.Range("M121").Value = -653
$ws.Range("N121").Value = -9893

$ws.Range("H129").Value = 1039.1562
$ws.Range("I129").Value = 296.33334
$ws.Range("J129").Value = 1210.5769
$ws.Range("K129").Value = 889.0000200000001
$ws.Range("L129").Value = 3631.7307
$ws.Range("M129").Value = 4110.99998
$ws.Range("N129").Value = -13631.7307

$ws.Range("H133").Value = 41272.727
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 41272.727
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 41272.727
$ws.Range("N133").Value = -51392.727

$ws.Range("H141").Value = 2553.0908
$ws.Range("I141").Value = 2113.5789
$ws.Range("J141").Value = 5336.6665
$ws.Range("K141").Value = 6340.736699999999
$ws.Range("L141").Value = 16009.9995
$ws.Range("M141").Value = -1160.736699999999
$ws.Range("N141").Value = -26369.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1992744.4
$ws.Range("I32").Value = 3763.018
$ws.Range("J32").Value = 15666991
$ws.Range("K32").Value = 3763.018
$ws.Range("L32").Value = 15666991
$ws.Range("M32").Value = -3476.018
$ws.Range("N32").Value = -15667565

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 466.33334
$ws.Range("I22").Value = 466.33334
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 466.33334
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -293.33334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 6727.727
$ws.Range("I62").Value = 6502.5
$ws.Range("J62").Value = 6777.778
$ws.Range("K62").Value = 6502.5
$ws.Range("L62").Value = 6777.778
$ws.Range("M62").Value = -5878.5
$ws.Range("N62").Value = -8025.778

$ws.Range("H65").Value = 6727.727
$ws.Range("I65").Value = 6502.5
$ws.Range("J65").Value = 6777.778
$ws.Range("K65").Value = 32512.5
$ws.Range("L65").Value = 33888.89
$ws.Range("M65").Value = -29392.5
$ws.Range("N65").Value = -40128.89

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 1971
$ws.Range("I64").Value = 1006
$ws.Range("J64").Value = 2453.5
$ws.Range("K64").Value = 3018
$ws.Range("L64").Value = 7360.5
$ws.Range("M64").Value = -2748
$ws.Range("N64").Value = -7900.5

$ws.Range("H67").Value = 1971
$ws.Range("I67").Value = 1006
$ws.Range("J67").Value = 2453.5
$ws.Range("K67").Value = 3018
$ws.Range("L67").Value = 7360.5
$ws.Range("M67").Value = -2082
$ws.Range("N67").Value = -9232.5

$ws.Range("H114").Value = 1563.2593
$ws.Range("I114").Value = 845
$ws.Range("J114").Value = 1922.3889
$ws.Range("K114").Value = 2535
$ws.Range("L114").Value = 5767.1667
$ws.Range("M114").Value = 719
$ws.Range("N114").Value = -12275.1667

$ws.Range("H139").Value = 1981.1177
$ws.Range("I139").Value = 1604.9375
$ws.Range("J139").Value = 8000
$ws.Range("K139").Value = 4814.8125
$ws.Range("L139").Value = 24000
$ws.Range("M139").Value = 325.1875
$ws.Range("N139").Value = -34280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 12052.714
$ws.Range("I57").Value = 6053.5
$ws.Range("J57").Value = 12684.211
$ws.Range("K57").Value = 6053.5
$ws.Range("L57").Value = 12684.211
$ws.Range("M57").Value = -5233.5
$ws.Range("N57").Value = -14324.211

$ws.Range("H70").Value = 4557.25
$ws.Range("I70").Value = 4092.2307
$ws.Range("J70").Value = 5420.857
$ws.Range("K70").Value = 4092.2307
$ws.Range("L70").Value = 5420.857
$ws.Range("M70").Value = -3822.2307
$ws.Range("N70").Value = -5960.857

$ws.Range("H73").Value = 4557.25
$ws.Range("I73").Value = 4092.2307
$ws.Range("J73").Value = 5420.857
$ws.Range("K73").Value = 4092.2307
$ws.Range("L73").Value = 5420.857
$ws.Range("M73").Value = -3156.2307
$ws.Range("N73").Value = -7292.857

$ws.Range("H80").Value = 2963.9583
$ws.Range("I80").Value = 3061.75
$ws.Range("J80").Value = 2475
$ws.Range("K80").Value = 3061.75
$ws.Range("L80").Value = 2475
$ws.Range("M80").Value = -2063.75
$ws.Range("N80").Value = -4471

$ws.Range("H83").Value = 2963.9583
$ws.Range("I83").Value = 3061.75
$ws.Range("J83").Value = 2475
$ws.Range("K83").Value = 15308.75
$ws.Range("L83").Value = 12375
$ws.Range("M83").Value = -10316.75
$ws.Range("N83").Value = -22359

$ws.Range("H86").Value = 8000
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 8000
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 8000
$ws.Range("N86").Value = -10372

$ws.Range("H89").Value = 8000
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 8000
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 24000
$ws.Range("N89").Value = -35856

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3150.3333
$ws.Range("I82").Value = 2850.3333
$ws.Range("J82").Value = 3750.3333
$ws.Range("K82").Value = 2850.3333
$ws.Range("L82").Value = 3750.3333
$ws.Range("M82").Value = -2489.3333
$ws.Range("N82").Value = -4472.3333

$ws.Range("H85").Value = 3150.3333
$ws.Range("I85").Value = 2850.3333
$ws.Range("J85").Value = 3750.3333
$ws.Range("K85").Value = 2850.3333
$ws.Range("L85").Value = 3750.3333
$ws.Range("M85").Value = -1602.3333
$ws.Range("N85").Value = -6246.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2333.3333
$ws.Range("I81").Value = 2000
$ws.Range("J81").Value = 3000
$ws.Range("K81").Value = 4000
$ws.Range("L81").Value = 6000
$ws.Range("M81").Value = -2939
$ws.Range("N81").Value = -8122

$ws.Range("H84").Value = 2333.3333
$ws.Range("I84").Value = 2000
$ws.Range("J84").Value = 3000
$ws.Range("K84").Value = 20000
$ws.Range("L84").Value = 30000
$ws.Range("M84").Value = -14696
$ws.Range("N84").Value = -40608

